$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Values are entered in this specific order so the shared-strings table is
# built up the same way as in the target workbook (column-major per
# semantic group, matching how the sheet was originally authored).

# StructureRelation column (S)
$ws.Range("S1").Value = "StructureRelation"
$ws.Range("S2").Value = "relation-1"
$ws.Range("S3").Value = "relation-2"
$ws.Range("S4").Value = "relation-3"

# Content / Format headers (T1, U1)
$ws.Range("T1").Value = "Content"
$ws.Range("U1").Value = "Format"

# Content values (T2:T4)
$ws.Range("T2").Value = "content-1"
$ws.Range("T3").Value = "content-2"
$ws.Range("T4").Value = "content-3"

# Format values (U2:U4)
$ws.Range("U2").Value = "format-1"
$ws.Range("U3").Value = "format-2"
$ws.Range("U4").Value = "format-3"

# Formula column (V)
$ws.Range("V1").Value = "Formula"
$ws.Range("V2").Value = "formula-1"
$ws.Range("V3").Value = "formula-2"
$ws.Range("V4").Value = "formula-3"

# Smiles column (W)
$ws.Range("W1").Value = "Smiles"
$ws.Range("W2").Value = "CCC-1"
$ws.Range("W3").Value = "CCC-2"
$ws.Range("W4").Value = "CCC-3"

# InChI column (X)
$ws.Range("X1").Value = "InChI"
$ws.Range("X2").Value = "inchi-1"
$ws.Range("X3").Value = "inchi-2"
$ws.Range("X4").Value = "inchi-3"

# InchiKey column (Y)
$ws.Range("Y1").Value = "InchiKey"
$ws.Range("Y2").Value = "inchi-key-1"
$ws.Range("Y3").Value = "inchi-key-2"
$ws.Range("Y4").Value = "inchi-key-3"

# property1 / property2 headers (Z1, AA1)
$ws.Range("Z1").Value = "property1"
$ws.Range("AA1").Value = "property2"

# property1 values (Z2:Z4)
$ws.Range("Z2").Value = "prop1-1"
$ws.Range("Z3").Value = "prop1-2"
$ws.Range("Z4").Value = "prop1-3"

# property3 header (AB1)
$ws.Range("AB1").Value = "property3"

# property2 numeric values (AA2:AA4)
$ws.Range("AA2").Value = 5.0999999999999996
$ws.Range("AA3").Value = 5.2
$ws.Range("AA4").Value = 5.3

# property3 numeric values (AB2:AB4)
$ws.Range("AB2").Value = 200.1
$ws.Range("AB3").Value = 200.2
$ws.Range("AB4").Value = 200.3

# Match formatting of the new header cells to the existing header style (R1 is bold, "s=1")
$ws.Range("R1").Copy()
$ws.Range("S1:AB1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the sheet view: scroll so column P is the top-left visible column,
# and select AB7 as the active cell.
$ws.Application.ActiveWindow.ScrollColumn = 16
$ws.Range("AB7").Select()
